# commit after login and batch
# Adds a new "login" worksheet (with username/password test fixtures and
# hyperlinked email addresses) after the existing "staff" sheet, and makes
# it the active/selected sheet.

$wb = $excel.ActiveWorkbook

# Add the new sheet after the last existing sheet ("staff") so it lands at
# the end of the tab strip.
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "login"

# Populate row 2 first (bottom-up, right-to-left) so the shared-string
# table is built in the same order as the source edit, then row 1.
$ws.Range("D2").Value = "mytesting"
$ws.Range("C2").Value = "bhuvee@mail.com"
$ws.Range("B2").Value = "UIHackathon@02"
$ws.Range("A2").Value = "sdetorganizers@gmail.com"

$ws.Range("D1").Value = "invalidPassword"
$ws.Range("C1").Value = "invalidUsername"
$ws.Range("B1").Value = "validPassword"
$ws.Range("A1").Value = "validUserName"

# Turn the three email cells in row 2 into mailto hyperlinks (this also
# applies the built-in "Hyperlink" style/font to them).
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:sdetorganizers@gmail.com", "", "", "sdetorganizers@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:UIHackathon@02", "", "", "UIHackathon@02")
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:bhuvee@mail.com", "", "", "bhuvee@mail.com")

# Make "login" the active sheet/tab and select D25 on it, matching the
# recorded selection state.
$ws.Activate()
$null = $ws.Range("D25").Select()

Write-Host "login sheet added"
